$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate cell values ---
# The fill order below mirrors the order new strings first appear in the
# target workbook (columns A-D mostly first, with a couple of deferred
# cells, then columns E-F) so that brand-new shared strings end up appended
# in the same relative order as the target file. (Strings that already
# existed in the workbook - Read-heavy, Write-heavy, Aggregate data, Mixed -
# are kept/deduplicated automatically by the engine.)

$ws.Range("A1").Value = "Entities"
$ws.Range("B1").Value = "Operations"
$ws.Range("C1").Value = "Information Needed"
$ws.Range("D1").Value = "Type"

$ws.Range("A3").Value = "Products"
$ws.Range("B3").Value = "Find by field"
$ws.Range("C3").Value = "Field: name"
$ws.Range("D3").Value = "Read-heavy"

$ws.Range("A4").Value = "Inventory"
$ws.Range("B4").Value = "Insert new document"
$ws.Range("D4").Value = "Write-heavy"

$ws.Range("A5").Value = "ProductionLines"
$ws.Range("B5").Value = "Update document"
$ws.Range("D5").Value = "Write-heavy"

$ws.Range("A6").Value = "ProductionTasks"
$ws.Range("B6").Value = "Aggregate data"
$ws.Range("C6").Value = "Pipeline stages for aggregation"
$ws.Range("D6").Value = "Mixed"

$ws.Range("C4").Value = "Fields: name, category, unit_price"
$ws.Range("C5").Value = "Filter criteria: _id, Field to update: status"

$ws.Range("B2").Value = "Find by objectID"
$ws.Range("C2").Value = "Field: _id"

$ws.Range("E1").Value = "Indexed Fields"
$ws.Range("E2").Value = "_id"
$ws.Range("E3").Value = "name"

$ws.Range("F2").Value = "Single Field Index"
$ws.Range("F1").Value = "Index Type"
$ws.Range("F3").Value = "Single Field Index"
$ws.Range("F4").Value = "No Index Needed"
$ws.Range("F5").Value = "Single Field Index"
$ws.Range("F6").Value = "No Index Needed"

$ws.Range("A2").Value = "Inventory"
$ws.Range("D2").Value = "Read-heavy"
$ws.Range("E5").Value = "_id"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 15.333333333333334
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws.Columns.Item(3).ColumnWidth = 34.666666666666664
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(6).ColumnWidth = 20.333333333333332

# --- View: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("G6").Select()
